$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.392.20"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "2.045.26"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.87%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.15"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.386"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0805"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.98%  "
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("D12").Value = "2.351.57"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.753"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("D17").Value = "2.046.27"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").Value = "37.292.12"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("D21").Value = "0.0₃0852"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("E25").Value = "  -4.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.127"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.99%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.26%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.60%  "
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.73%  "
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.25%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.30%  "
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("E40").Value = "  -7.31%  "
$ws.Range("D41").Value = "1.499.13"
$ws.Range("E41").Value = "  +3.18%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0938"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.48%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.68%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.80%  "
$ws.Range("D51").Value = "2.239.41"
$ws.Range("E51").Value = "  -1.86%  "
